# feat: add 2022-Q1 data
#
# Layout before:  [2020-Q4] [总计]
# Layout after:   [2020-Q4] [2022-Q1] [总计]
#
# The existing "总计" sheet (position 2) becomes the new "2022-Q1"
# fund-holdings detail sheet (it keeps the sheetId/rId it already had), and
# a clone of the original "总计" sheet is appended at the end, renamed back
# to "总计", with the refreshed totals table (2022-Q1 row added on top,
# 2020-Q4 row pushed down to row 3).

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Grab format-donor cells from the current "总计" sheet BEFORE touching
# any values - the *style index* they carry (bold header / bordered
# index cell) is what we want to replicate onto the new cells. Content
# changes later don't affect a style index already handed out.
# ---------------------------------------------------------------------
$headerStyleSrc = $zongji.Range("B1")
$indexStyleSrc  = $zongji.Range("A2")

# ===========================================================================
# STEP 1 - clone the current "总计" sheet to the end of the workbook. This
# clone keeps all of "总计"'s boilerplate (sheetPr/pageMargins/styles) and
# its existing 2020-Q4 row, which is exactly what the refreshed "总计" sheet
# needs as a starting point. Give it a placeholder name for now so it
# doesn't collide with the source sheet (which is still called "总计").
# ===========================================================================
$zongji.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$zj = $wb.Worksheets.Item($wb.Worksheets.Count)
$zj.Name = "总计新"

# ===========================================================================
# STEP 2 - repurpose the ORIGINAL "总计" sheet into "2022-Q1"
# ===========================================================================
$q1 = $zongji
$q1.Name = "2022-Q1"

# ---- header row ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerStyleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# ---- data row (row 2) ----
$indexStyleSrc.Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A2").Value = 0

$textRow = $q1.Range("B2:G2")
$textRow.NumberFormat = "@"
$q1.Range("B2").Value = "160324"
$q1.Range("C2").Value = "华夏磐晟灵活配置混合（LOF）"
$q1.Range("D2").Value = "1.07"
$q1.Range("E2").Value = "78.29"
$q1.Range("F2").Value = "3.92"
$q1.Range("G2").Value = "0.0419"
$textRow.ClearFormats()

$q1.Range("H2").Value = 9

$excel.CutCopyMode = 0

# ===========================================================================
# STEP 3 - finish turning the clone into the refreshed "总计" sheet: push the
# existing 2020-Q4 totals row down to row 3, then write the new 2022-Q1
# totals into row 2. (Row 3's values are the clone's pre-existing 2020-Q4
# totals, written out explicitly rather than round-tripped through a read,
# since they already occupy row 2 of the clone at this point.)
# ===========================================================================
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.02

$zj.Range("B3").NumberFormat = "@"
$zj.Range("B3").Value = "2020-Q4"
$zj.Range("B3").ClearFormats()

$indexStyleSrc.Copy()
$zj.Range("A3").PasteSpecial(-4122)
$zj.Range("A3").Value = 1

$indexStyleSrc.Copy()
$zj.Range("A2").PasteSpecial(-4122)
$zj.Range("A2").Value = 0

$zj.Range("B2").NumberFormat = "@"
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("B2").ClearFormats()
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.04

$zj.Name = "总计"

$excel.CutCopyMode = 0

# ===========================================================================
# STEP 4 - restore the original active sheet/selection (unchanged by this
# edit: "2020-Q4" stays the active tab).
# ===========================================================================
$wb.Worksheets.Item(1).Activate()
